$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 12 ("Accuracy Report"): nudge the KNN-accuracy group shape right.
#    Target OOXML: <a:off x="4447439" y="2206807"/> (was x="4424814").
#    Shape.Left is in points; 4447439 EMU / 12700 = 350.19204724... pt, but
#    PowerPoint's point->EMU rounding needs a hair more precision to land on
#    the exact target EMU value, so nudge slightly above the midpoint.
# ---------------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$grp = $slide12.Shapes.Item(3)
$grp.Left = 350.19207

# ---------------------------------------------------------------------------
# 2) Slide 3 ("Darsh #1 Fan Club!  XOXO"): retitle.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "Darsh’s only fans!  XOXO"

# ---------------------------------------------------------------------------
# 3) Slide 9 ("Modeling"): hide the slide from the slideshow.
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$slide9.SlideShowTransition.Hidden = $true

# ---------------------------------------------------------------------------
# 4) Slide 9: prefix the "Models tried" label with four spaces (new run).
# ---------------------------------------------------------------------------
$label = $slide9.Shapes.Item(2)
$label.TextFrame.TextRange.InsertBefore("    ")

# ---------------------------------------------------------------------------
# 5) Slide 9: swap the table's style to the new tableStyleId.
# ---------------------------------------------------------------------------
$tableShape = $slide9.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{FD569BA4-3733-4328-87F6-26B5889362EF}")

# ---------------------------------------------------------------------------
# 6) Theme: the deck's active theme (master/theme2.xml) becomes the plain
#    "Default" color scheme instead of the Office palette.
# ---------------------------------------------------------------------------
$tcs = $slide3.ThemeColorScheme
$tcs.Item(1).RGB = 0        # dk1      = 000000
$tcs.Item(2).RGB = 16777215 # lt1      = FFFFFF
$tcs.Item(3).RGB = 5800213  # dk2      = 158158
$tcs.Item(4).RGB = 15987699 # lt2      = F3F3F3
$tcs.Item(5).RGB = 13077765 # accent1  = 058DC7
$tcs.Item(6).RGB = 3322960  # accent2  = 50B432
$tcs.Item(7).RGB = 1791725  # accent3  = ED561B
$tcs.Item(8).RGB = 61421    # accent4  = EDEF00
$tcs.Item(9).RGB = 15059748 # accent5  = 24CBE5
$tcs.Item(10).RGB = 7529828 # accent6  = 64E572
$tcs.Item(11).RGB = 13369378 # hlink    = 2200CC
$tcs.Item(12).RGB = 9116245  # folHlink = 551A8B
